# Refresh the cryptos price/volume table to the latest snapshot
# (coin rows 13/14 and 39/40 also swap rank order).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.249.31"
$ws.Range("E2").Value = "  -2.56%  "

$ws.Range("D3").Value = "1.559.21"
$ws.Range("E3").Value = "  -3.94%  "

$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'206.54"
$ws.Range("E5").Value = "  -3.27%  "

$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").Value = "'0.477"
$ws.Range("E7").Value = "  -5.08%  "

$ws.Range("D8").Value = "'0.0604"
$ws.Range("E8").Value = "  -1.65%  "

$ws.Range("E9").Value = "  -3.38%  "

$ws.Range("D10").Value = "'17.70"
$ws.Range("E10").Value = "  -3.06%  "

$ws.Range("D11").Value = "'0.0781"
$ws.Range("E11").Value = "  -0.76%  "

$ws.Range("D12").Value = "1.780.39"
$ws.Range("E12").Value = "  -3.69%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.563.04"
$ws.Range("E13").Value = "  -3.69%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'3.99"
$ws.Range("E14").Value = "  -4.44%  "

$ws.Range("D15").Value = "'0.502"
$ws.Range("E15").Value = "  -4.06%  "

$ws.Range("D16").Value = "25.253.02"
$ws.Range("E16").Value = "  -2.53%  "

$ws.Range("D17").Value = "'58.96"
$ws.Range("E17").Value = "  -3.51%  "

$ws.Range("D18").Value = "0.0₃0708"
$ws.Range("E18").Value = "  -3.29%  "

$ws.Range("E19").Value = "  -0.33%  "

$ws.Range("D20").Value = "'184.98"
$ws.Range("E20").Value = "  -3.54%  "

$ws.Range("D21").Value = "'4.10"
$ws.Range("E21").Value = "  -3.08%  "

$ws.Range("D22").Value = "'9.25"
$ws.Range("E22").Value = "  -3.27%  "

$ws.Range("D23").Value = "'5.85"
$ws.Range("E23").Value = "  -3.51%  "

$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("E25").Value = "  -3.86%  "

$ws.Range("D26").Value = "'139.50"
$ws.Range("E26").Value = "  -3.00%  "

$ws.Range("D27").Value = "'1.64"
$ws.Range("E27").Value = "  -5.99%  "

$ws.Range("D28").Value = "'6.42"
$ws.Range("E28").Value = "  -4.44%  "

$ws.Range("D29").Value = "'14.78"
$ws.Range("E29").Value = "  -2.27%  "

$ws.Range("E30").Value = "  -6.44%  "

$ws.Range("D31").Value = "'0.0463"
$ws.Range("E31").Value = "  -3.90%  "

$ws.Range("D32").Value = "'3.02"
$ws.Range("E32").Value = "  -2.97%  "

$ws.Range("D33").Value = "'2.97"
$ws.Range("E33").Value = "  -4.20%  "

$ws.Range("D34").Value = "'1.45"
$ws.Range("E34").Value = "  -2.62%  "

$ws.Range("E35").Value = "  -4.02%  "

$ws.Range("D36").Value = "1.086.26"
$ws.Range("E36").Value = "  -2.69%  "

$ws.Range("E38").Value = "  -5.01%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'0.826"
$ws.Range("E39").Value = "  +7.82%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0149"
$ws.Range("E40").Value = "  -2.23%  "

$ws.Range("E41").Value = "  -4.60%  "

$ws.Range("D42").Value = "'0.760"
$ws.Range("E42").Value = "  -10.09%  "

$ws.Range("D43").Value = "'92.76"
$ws.Range("E43").Value = "  -5.28%  "

$ws.Range("D45").Value = "1.694.48"
$ws.Range("E45").Value = "  -3.68%  "

$ws.Range("E46").Value = "  -2.70%  "

$ws.Range("D47").Value = "'52.33"
$ws.Range("E47").Value = "  -3.64%  "

$ws.Range("D48").Value = "'0.0504"
$ws.Range("E48").Value = "  -4.94%  "

$ws.Range("E49").Value = "  -2.13%  "

$ws.Range("D50").Value = "'0.405"
$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("E51").Value = "  -0.42%  "
